$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contrib")

$ws.Range("B2").Value = 24.74752479652351
$ws.Range("B3").Value = 17.54282074558445
$ws.Range("B4").Value = 16.83120140573486
$ws.Range("B5").Value = 12.90192333711556
$ws.Range("B6").Value = 11.15349735918752
$ws.Range("B7").Value = 7.166798298852618
$ws.Range("B8").Value = 5.857640853090421
$ws.Range("B9").Value = 3.798593203911072
